$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12486
$ws1.Range("F3").Value = 6992
$ws1.Range("F9").Value = 15
$ws1.Range("F10").Value = 973
$ws1.Range("F11").Value = 126
$ws1.Range("F12").Value = 330
$ws1.Range("F15").Value = 1003
$ws1.Range("F17").Value = 222
$ws1.Range("F20").Value = 261
$ws1.Range("F22").Value = 36
$ws1.Range("F23").Value = 102
$ws1.Range("F25").Value = 5150
$ws1.Range("F26").Value = 64
$ws1.Range("F27").Value = 1375
$ws1.Range("F28").Value = 283
$ws1.Range("F29").Value = 1173
$ws1.Range("G29").Value = 67.5
$ws1.Range("F30").Value = 1301
$ws1.Range("F31").Value = 577
$ws1.Range("F33").Value = 3711

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 3736
$ws2.Range("G4").Value = "不可售"
$ws2.Range("F6").Value = 14
$ws2.Range("F7").Value = 28

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9212
$ws3.Range("F4").Value = 1941

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9212
$ws4.Range("F4").Value = 1941
$ws4.Range("F5").Value = 12487
$ws4.Range("F6").Value = 6992
$ws4.Range("F8").Value = 3736
$ws4.Range("G8").Value = "不可售"
$ws4.Range("F14").Value = 15
$ws4.Range("F15").Value = 973
$ws4.Range("F16").Value = 126
$ws4.Range("F17").Value = 330
$ws4.Range("F20").Value = 1003
$ws4.Range("F22").Value = 222
$ws4.Range("F25").Value = 261
$ws4.Range("F27").Value = 36
$ws4.Range("F29").Value = 14
$ws4.Range("F33").Value = 64
$ws4.Range("F34").Value = 1375
$ws4.Range("F37").Value = 283
$ws4.Range("F39").Value = 1173
$ws4.Range("G39").Value = 67.5
$ws4.Range("F40").Value = 1301
$ws4.Range("F41").Value = 577
$ws4.Range("F46").Value = 3711
